# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2
# of the zh-cn and de-de sheets, as part of generating the
# handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 06:51:43"
$wsZhCn.Range("H2").Value = "2016-03-20 06:52:01"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 06:51:46"
$wsDeDe.Range("H2").Value = "2016-03-20 06:52:07"
